$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New report data (rows 2-7). Column C (dates) must stay plain text, so
# force the number format to Text before writing those values.
$data = @(
    @("cli-9876", "ana schmidt", "10/10/2025", "r$ 250,75", "['- item a: 1 unidade', '- item b: 3 unidades']"),
    @("cli-9876", "ana schmidt", "10/10/2025", "r$ 250,75", "['- item a: 1 unidade', '- item b: 3 unidades']"),
    @("cli-5432", "bruno costa", "12/10/2025", "r$ 1.120,00", "['- item c: 10 unidades', '- item d: 1 unidade']"),
    @("cli-9876", "ana schmidt", "10/10/2025", "r$ 250,75", "['- item a: 1 unidade', '- item b: 3 unidades']"),
    @("cli-9876", "ana schmidt", "10/10/2025", "r$ 250,75", "['- item a: 1 unidade', '- item b: 3 unidades']"),
    @("cli-5432", "bruno costa", "12/10/2025", "r$ 1.120,00", "['- item c: 10 unidades', '- item d: 1 unidade']")
)

$ws.Range("C2:C7").NumberFormat = "@"

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $row++
}

# Drop the temporary Text number format again so the cells end up
# unstyled (matching the source rows, which carry no explicit style).
$ws.Range("C2:C7").Style = "Normal"
